# Update strategic authority scores for the "Greater London Authority" row (row 3)
# as described in the commit message / diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 0.4761904761904762
$ws.Range("H3").Value = 0.8888888888888888
$ws.Range("I3").Value = 0.8
$ws.Range("N3").Value = 0.6579761904761904
